{"js": "// The document had a trailing empty paragraph right after the paragraph\n// that ends with \"...stay connected.\" and an orphaned \"_GoBack\" bookmark\n// sitting in the very last paragraph of the document. The edit moves the\n// \"_GoBack\" bookmark so it wraps the end of the \"...stay connected.\"\n// paragraph (this is what Word does when the cursor was last left there),\n// and removes the now-redundant empty paragraph that used to follow it.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst items = paragraphs.items;\nfor (const p of items) {\n  p.load(\"text\");\n}\nawait context.sync();\n\n// Locate the paragraph that ends the sentence about staying connected.\nlet targetIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.indexOf(\"stay connected.\") !== -1) {\n    targetIndex = i;\n    break;\n  }\n}\n\nif (targetIndex === -1) {\n  throw new Error(\"Could not locate the 'stay connected.' paragraph.\");\n}\n\n// Remove the existing \"_GoBack\" bookmark from wherever it currently lives\n// (the last paragraph in the document, trailing a tab character).\ncontext.document.deleteBookmark(\"_GoBack\");\n\n// Delete the empty paragraph that immediately follows the target\n// paragraph - it no longer carries the bookmark and is redundant.\nconst emptyParagraph = items[targetIndex + 1];\nemptyParagraph.delete();\n\n// Re-insert \"_GoBack\" as a collapsed bookmark at the end of the target\n// paragraph's text (before its paragraph mark).\nconst endRange = items[targetIndex].getRange(\"End\");\nendRange.insertBookmark(\"_GoBack\");\n\nawait context.sync();\n", "ps1": "# The document had a trailing empty paragraph right after the paragraph\n# that ends with \"...stay connected.\" and an orphaned \"_GoBack\" bookmark\n# sitting in the very last paragraph of the document. This edit moves the\n# \"_GoBack\" bookmark so it wraps the end of the \"...stay connected.\"\n# paragraph (mirroring where Word leaves it after the last edit), and\n# removes the now-redundant empty paragraph that used to follow it.\n\n$d = $word.ActiveDocument\n\n# Step 1: locate the paragraph that ends with the \"stay connected.\" sentence.\n$targetIndex = 0\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text -like \"*stay connected.*\") {\n        $targetIndex = $i\n        break\n    }\n}\nif ($targetIndex -eq 0) {\n    throw \"Could not locate the 'stay connected.' paragraph.\"\n}\n\n# Step 2: remove the existing \"_GoBack\" bookmark from wherever it currently\n# lives (the last paragraph in the document, trailing a tab character).\ntry {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n} catch {\n    # No pre-existing bookmark - nothing to remove.\n}\n\n# Step 3: delete the empty paragraph that immediately follows the target\n# paragraph - it no longer carries the bookmark and is redundant.\n$nextPara = $d.Paragraphs.Item($targetIndex + 1)\n$nextPara.Range.Delete()\n\n# Step 4: re-insert \"_GoBack\" as a collapsed bookmark at the end of the\n# target paragraph's text (right before its paragraph mark).\n#\n# A collapsed (zero-length) Range that sits exactly on the last character\n# slot of a paragraph confuses Bookmarks.Add in this host, so a single\n# placeholder character is temporarily appended, the bookmark is added in\n# front of it, and the placeholder is removed again.\n$targetPara = $d.Paragraphs.Item($targetIndex)\n$paraRange = $targetPara.Range\n$padRange = $d.Range($paraRange.End - 1, $paraRange.End - 1)\n$padRange.InsertAfter(\"X\")\n\n$bookmarkPos = $d.Range($paraRange.End - 2, $paraRange.End - 2)\n$d.Bookmarks.Add(\"_GoBack\", $bookmarkPos)\n\n$cleanupRange = $d.Range($paraRange.End - 2, $paraRange.End - 1)\n$cleanupRange.Delete()\n"}
